# Update the CDA Logical model for ST.r2b
# - Bump Version and Date metadata values
# - Insert a new "Jurisdiction" property row into the Metadata sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Version (row 3) and Date (row 8) values
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row above row 11 ("Description") for the new "Jurisdiction" property
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Re-apply the standard data-row formatting (border/alignment) to the new row,
# matching the rest of the property table.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

Write-Host "edit applied"
